# Daily attendance processing - 2026-01-14 11:08:54
# Swap the "Recorded By" display order in column G from
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# across every row of the "Session Analysis Results" sheet that has that
# exact value (rows where both System and the user recorded attendance).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
